{"js": "// Replace the two-digit multiplication problems in the worksheet table.\n// The document has a single table; every 5th row (0, 4, 9, 14, 19) holds\n// the five math expressions for that row, the rows in between are blank\n// (left for the student's handwritten answer). Each expression cell's\n// text is replaced in place with its new value, matching the order the\n// cells appear in the document.\n\nconst newValuesByRow = {\n  0: [\"64\u00d761=\", \"31\u00d728=\", \"70\u00d772=\", \"50\u00d798=\", \"88\u00d791=\"],\n  4: [\"97\u00d740=\", \"92\u00d768=\", \"69\u00d739=\", \"55\u00d785=\", \"56\u00d787=\"],\n  9: [\"20\u00d797=\", \"83\u00d783=\", \"80\u00d718=\", \"85\u00d745=\", \"75\u00d722=\"],\n  14: [\"73\u00d784=\", \"74\u00d748=\", \"18\u00d775=\", \"65\u00d757=\", \"43\u00d742=\"],\n  19: [\"55\u00d785=\", \"50\u00d794=\", \"12\u00d726=\", \"84\u00d742=\", \"40\u00d724=\"],\n};\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const rowIndex of Object.keys(newValuesByRow)) {\n  const idx = Number(rowIndex);\n  const newValues = newValuesByRow[idx];\n  const cells = rows.items[idx].cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (let col = 0; col < newValues.length; col++) {\n    cells.items[col].value = newValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the two-digit multiplication problems in the worksheet table.\n# The document has a single table; every 5th row (1, 5, 10, 15, 20 in\n# Word's 1-based Table.Cell indexing) holds the five math expressions for\n# that row, the rows in between are blank (left for the student's\n# handwritten answer). Each expression cell's text is replaced in place\n# with its new value, addressed directly by (row, column) so the two\n# duplicate \"18\u00d765=\" cells land on their own distinct replacements.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValuesByRow = @{\n    1  = @(\"64\u00d761=\", \"31\u00d728=\", \"70\u00d772=\", \"50\u00d798=\", \"88\u00d791=\")\n    5  = @(\"97\u00d740=\", \"92\u00d768=\", \"69\u00d739=\", \"55\u00d785=\", \"56\u00d787=\")\n    10 = @(\"20\u00d797=\", \"83\u00d783=\", \"80\u00d718=\", \"85\u00d745=\", \"75\u00d722=\")\n    15 = @(\"73\u00d784=\", \"74\u00d748=\", \"18\u00d775=\", \"65\u00d757=\", \"43\u00d742=\")\n    20 = @(\"55\u00d785=\", \"50\u00d794=\", \"12\u00d726=\", \"84\u00d742=\", \"40\u00d724=\")\n}\n\nforeach ($row in $newValuesByRow.Keys) {\n    $values = $newValuesByRow[$row]\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $t.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
